$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at the top; this shifts the existing header row
# (row 1 -> row 2) and all data rows (rows 2-27 -> rows 3-28) down by one,
# carrying their values AND their formatting along with them.
$ws.Rows.Item(1).Insert()

# The header formatting (bold / centered / bordered, style index 1) is now
# sitting on row 2 (it moved down with the old header row). Copy that
# formatting up onto the brand-new row 1 before we overwrite its values.
$ws.Range("A2:L2").Copy()
$ws.Range("A1:L1").PasteSpecial(-4122)

# Populate the new row 1 with the numeric column-index values 0..11.
for ($i = 0; $i -lt 12; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $i
}

# Row 2 (the old header row) should no longer carry the special header
# formatting - it reverts to the default/normal formatting.
$ws.Range("A2:L2").ClearFormats()

# In the old header row, the "thread_size" / "material_surface" labels in
# K/L are removed (left blank), matching the existing blank I column.
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
